# Weekly roll of the "Perejil" price series.
#
# The sheet holds one row-pair (Primera / Segunda calidad) per observation
# date, ordered from most-recent (top, row 2) to oldest (bottom, row 89).
# A new week's observation is inserted at the top of the historical window
# (rows 14-15, right after the 6 most-recent weeks that are left untouched),
# every older pair slides down by one pair (2 rows), and the oldest pair
# that falls off the bottom (old rows 88-89) is appended as brand-new rows
# 90-91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The pair that falls off the bottom of the existing window becomes two
#    new rows appended at the end of the sheet (old A88:R89 -> A90:R91).
$ws.Range("A88:R89").Copy($ws.Range("A90:R91"))

# 2) Slide every pair from row 16 down to row 88 downward by one pair,
#    pulling from the pair that currently sits two rows above it. Walk from
#    the bottom up so a pair is never overwritten before it has been read.
for ($r = 88; $r -ge 16; $r -= 2) {
    $ws.Range(("A{0}:R{1}" -f ($r - 2), ($r - 1))).Copy($ws.Range(("A{0}:R{1}" -f $r, ($r + 1))))
}

# 3) The newest pair (rows 14-15) gets this week's date; every other field
#    on those rows is unchanged.
$ws.Range("D14:D15").Value2 = 44425
